$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Date" value ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value2 = "2025-07-25T07:22:51+00:00"

# --- Elements sheet: update canonical terminology URLs ---
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("Z3").Value2 = "https://mos.esante.gouv.fr/NOS/TRE_R14-TypeDiplome/FHIR/TRE-R14-TypeDiplome?vs"
$wsElem.Range("Z4").Value2 = "https://mos.esante.gouv.fr/NOS/TRE_R16-LieuFormation/FHIR/TRE-R16-LieuFormation?vs"
$wsElem.Range("Z7").Value2 = "https://mos.esante.gouv.fr/NOS/TRE_R57-DiplomeEuropeenEtudeSpecialisee/FHIR/TRE-R57-DiplomeEuropeenEtudeSpecialisee?vs"

# The longer URL text widens column Z's stored "best fit" width
# (83.69140625 -> 103.24609375 chars). Re-apply the column width so the
# stored width reflects the new content length.
$wsElem.Columns.Item(26).ColumnWidth = 102.333333333
